# Daily Update Feb 08 2020 - add new sheet Feb07_0813pm with latest COVID-19 case data
$wb = $excel.ActiveWorkbook

# 1) Duplicate the most-recent sheet (Feb06_0805pm) to seed the new snapshot
$srcSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $srcSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Feb07_0813pm"

# 2) Refresh every row with the Feb 7, 8:13pm case counts (adds Cruise Ship / Others row)
$newSheet.Cells.Item(1,1).Value = 'State'
$newSheet.Cells.Item(1,2).Value = 'Country'
$newSheet.Cells.Item(1,3).Value = 'Last Update'
$newSheet.Cells.Item(1,4).Value = 'Confirmed'
$newSheet.Cells.Item(1,5).Value = 'Deaths'
$newSheet.Cells.Item(1,6).Value = 'Recovered'
$newSheet.Cells.Item(2,1).Value = 'Hubei'
$newSheet.Cells.Item(2,2).Value = 'Mainland China'
$newSheet.Cells.Item(2,3).Value = 43868.58541666667
$newSheet.Cells.Item(2,4).Value = 22112
$newSheet.Cells.Item(2,5).Value = 618
$newSheet.Cells.Item(2,6).Value = 867
$newSheet.Cells.Item(3,1).Value = 'Guangdong'
$newSheet.Cells.Item(3,2).Value = 'Mainland China'
$newSheet.Cells.Item(3,3).Value = 43868.42569444444
$newSheet.Cells.Item(3,4).Value = 1034
$newSheet.Cells.Item(3,5).Value = 1
$newSheet.Cells.Item(3,6).Value = 88
$newSheet.Cells.Item(4,1).Value = 'Zhejiang'
$newSheet.Cells.Item(4,2).Value = 'Mainland China'
$newSheet.Cells.Item(4,3).Value = 43868.48125
$newSheet.Cells.Item(4,4).Value = 1006
$newSheet.Cells.Item(4,5).Value = 0
$newSheet.Cells.Item(4,6).Value = 123
$newSheet.Cells.Item(5,1).Value = 'Henan'
$newSheet.Cells.Item(5,2).Value = 'Mainland China'
$newSheet.Cells.Item(5,3).Value = 43868.58541666667
$newSheet.Cells.Item(5,4).Value = 914
$newSheet.Cells.Item(5,5).Value = 3
$newSheet.Cells.Item(5,6).Value = 86
$newSheet.Cells.Item(6,1).Value = 'Hunan'
$newSheet.Cells.Item(6,2).Value = 'Mainland China'
$newSheet.Cells.Item(6,3).Value = 43868.48125
$newSheet.Cells.Item(6,4).Value = 772
$newSheet.Cells.Item(6,5).Value = 0
$newSheet.Cells.Item(6,6).Value = 112
$newSheet.Cells.Item(7,1).Value = 'Anhui'
$newSheet.Cells.Item(7,2).Value = 'Mainland China'
$newSheet.Cells.Item(7,3).Value = 43868.217361111114
$newSheet.Cells.Item(7,4).Value = 665
$newSheet.Cells.Item(7,5).Value = 0
$newSheet.Cells.Item(7,6).Value = 47
$newSheet.Cells.Item(8,1).Value = 'Jiangxi'
$newSheet.Cells.Item(8,2).Value = 'Mainland China'
$newSheet.Cells.Item(8,3).Value = 43868.05763888889
$newSheet.Cells.Item(8,4).Value = 661
$newSheet.Cells.Item(8,5).Value = 0
$newSheet.Cells.Item(8,6).Value = 45
$newSheet.Cells.Item(9,1).Value = 'Chongqing'
$newSheet.Cells.Item(9,2).Value = 'Mainland China'
$newSheet.Cells.Item(9,3).Value = 43868.50208333333
$newSheet.Cells.Item(9,4).Value = 415
$newSheet.Cells.Item(9,5).Value = 2
$newSheet.Cells.Item(9,6).Value = 31
$newSheet.Cells.Item(10,1).Value = 'Jiangsu'
$newSheet.Cells.Item(10,2).Value = 'Mainland China'
$newSheet.Cells.Item(10,3).Value = 43868.58541666667
$newSheet.Cells.Item(10,4).Value = 408
$newSheet.Cells.Item(10,5).Value = 0
$newSheet.Cells.Item(10,6).Value = 43
$newSheet.Cells.Item(11,1).Value = 'Shandong'
$newSheet.Cells.Item(11,2).Value = 'Mainland China'
$newSheet.Cells.Item(11,3).Value = 43868.48125
$newSheet.Cells.Item(11,4).Value = 386
$newSheet.Cells.Item(11,5).Value = 0
$newSheet.Cells.Item(11,6).Value = 37
$newSheet.Cells.Item(12,1).Value = 'Sichuan'
$newSheet.Cells.Item(12,2).Value = 'Mainland China'
$newSheet.Cells.Item(12,3).Value = 43868.34930555556
$newSheet.Cells.Item(12,4).Value = 344
$newSheet.Cells.Item(12,5).Value = 1
$newSheet.Cells.Item(12,6).Value = 42
$newSheet.Cells.Item(13,1).Value = 'Beijing'
$newSheet.Cells.Item(13,2).Value = 'Mainland China'
$newSheet.Cells.Item(13,3).Value = 43868.14097222222
$newSheet.Cells.Item(13,4).Value = 297
$newSheet.Cells.Item(13,5).Value = 1
$newSheet.Cells.Item(13,6).Value = 33
$newSheet.Cells.Item(14,1).Value = 'Heilongjiang'
$newSheet.Cells.Item(14,2).Value = 'Mainland China'
$newSheet.Cells.Item(14,3).Value = 43868.32152777778
$newSheet.Cells.Item(14,4).Value = 277
$newSheet.Cells.Item(14,5).Value = 3
$newSheet.Cells.Item(14,6).Value = 12
$newSheet.Cells.Item(15,1).Value = 'Shanghai'
$newSheet.Cells.Item(15,2).Value = 'Mainland China'
$newSheet.Cells.Item(15,3).Value = 43868.259722222225
$newSheet.Cells.Item(15,4).Value = 277
$newSheet.Cells.Item(15,5).Value = 1
$newSheet.Cells.Item(15,6).Value = 30
$newSheet.Cells.Item(16,1).Value = 'Fujian'
$newSheet.Cells.Item(16,2).Value = 'Mainland China'
$newSheet.Cells.Item(16,3).Value = 43868.55763888889
$newSheet.Cells.Item(16,4).Value = 224
$newSheet.Cells.Item(16,5).Value = 0
$newSheet.Cells.Item(16,6).Value = 20
$newSheet.Cells.Item(17,1).Value = 'Shaanxi'
$newSheet.Cells.Item(17,2).Value = 'Mainland China'
$newSheet.Cells.Item(17,3).Value = 43868.34930555556
$newSheet.Cells.Item(17,4).Value = 184
$newSheet.Cells.Item(17,5).Value = 0
$newSheet.Cells.Item(17,6).Value = 17
$newSheet.Cells.Item(18,1).Value = 'Guangxi'
$newSheet.Cells.Item(18,2).Value = 'Mainland China'
$newSheet.Cells.Item(18,3).Value = 43868.02291666667
$newSheet.Cells.Item(18,4).Value = 172
$newSheet.Cells.Item(18,5).Value = 0
$newSheet.Cells.Item(18,6).Value = 17
$newSheet.Cells.Item(19,1).Value = 'Hebei'
$newSheet.Cells.Item(19,2).Value = 'Mainland China'
$newSheet.Cells.Item(19,3).Value = 43868.50208333333
$newSheet.Cells.Item(19,4).Value = 172
$newSheet.Cells.Item(19,5).Value = 1
$newSheet.Cells.Item(19,6).Value = 22
$newSheet.Cells.Item(20,1).Value = 'Yunnan'
$newSheet.Cells.Item(20,2).Value = 'Mainland China'
$newSheet.Cells.Item(20,3).Value = 43868.97430555556
$newSheet.Cells.Item(20,4).Value = 138
$newSheet.Cells.Item(20,5).Value = 0
$newSheet.Cells.Item(20,6).Value = 12
$newSheet.Cells.Item(21,1).Value = 'Hainan'
$newSheet.Cells.Item(21,2).Value = 'Mainland China'
$newSheet.Cells.Item(21,3).Value = 43868.52291666667
$newSheet.Cells.Item(21,4).Value = 117
$newSheet.Cells.Item(21,5).Value = 2
$newSheet.Cells.Item(21,6).Value = 10
$newSheet.Cells.Item(22,1).Value = 'Liaoning'
$newSheet.Cells.Item(22,2).Value = 'Mainland China'
$newSheet.Cells.Item(22,3).Value = 43868.64791666667
$newSheet.Cells.Item(22,4).Value = 99
$newSheet.Cells.Item(22,5).Value = 0
$newSheet.Cells.Item(22,6).Value = 7
$newSheet.Cells.Item(23,1).Value = 'Shanxi'
$newSheet.Cells.Item(23,2).Value = 'Mainland China'
$newSheet.Cells.Item(23,3).Value = 43868.52291666667
$newSheet.Cells.Item(23,4).Value = 96
$newSheet.Cells.Item(23,5).Value = 0
$newSheet.Cells.Item(23,6).Value = 15
$newSheet.Cells.Item(24,1).Value = 'Guizhou'
$newSheet.Cells.Item(24,2).Value = 'Mainland China'
$newSheet.Cells.Item(24,3).Value = 43868.467361111114
$newSheet.Cells.Item(24,4).Value = 81
$newSheet.Cells.Item(24,5).Value = 1
$newSheet.Cells.Item(24,6).Value = 6
$newSheet.Cells.Item(25,1).Value = 'Tianjin'
$newSheet.Cells.Item(25,2).Value = 'Mainland China'
$newSheet.Cells.Item(25,3).Value = 43868.23819444444
$newSheet.Cells.Item(25,4).Value = 81
$newSheet.Cells.Item(25,5).Value = 1
$newSheet.Cells.Item(25,6).Value = 2
$newSheet.Cells.Item(26,1).Value = 'Gansu'
$newSheet.Cells.Item(26,2).Value = 'Mainland China'
$newSheet.Cells.Item(26,3).Value = 43868.10625
$newSheet.Cells.Item(26,4).Value = 67
$newSheet.Cells.Item(26,5).Value = 0
$newSheet.Cells.Item(26,6).Value = 9
$newSheet.Cells.Item(27,1).Value = 'Jilin'
$newSheet.Cells.Item(27,2).Value = 'Mainland China'
$newSheet.Cells.Item(27,3).Value = 43868.029861111114
$newSheet.Cells.Item(27,4).Value = 65
$newSheet.Cells.Item(27,5).Value = 1
$newSheet.Cells.Item(27,6).Value = 4
$newSheet.Cells.Item(28,1).Value = 'Cruise Ship'
$newSheet.Cells.Item(28,2).Value = 'Others'
$newSheet.Cells.Item(28,3).Value = 43868.68958333333
$newSheet.Cells.Item(28,4).Value = 61
$newSheet.Cells.Item(28,5).Value = 0
$newSheet.Cells.Item(28,6).Value = 0
$newSheet.Cells.Item(29,1).Value = 'Inner Mongolia'
$newSheet.Cells.Item(29,2).Value = 'Mainland China'
$newSheet.Cells.Item(29,3).Value = 43868.467361111114
$newSheet.Cells.Item(29,4).Value = 50
$newSheet.Cells.Item(29,5).Value = 0
$newSheet.Cells.Item(29,6).Value = 5
$newSheet.Cells.Item(30,1).Value = 'Ningxia'
$newSheet.Cells.Item(30,2).Value = 'Mainland China'
$newSheet.Cells.Item(30,3).Value = 43868.467361111114
$newSheet.Cells.Item(30,4).Value = 43
$newSheet.Cells.Item(30,5).Value = 0
$newSheet.Cells.Item(30,6).Value = 5
$newSheet.Cells.Item(31,1).Value = 'Xinjiang'
$newSheet.Cells.Item(31,2).Value = 'Mainland China'
$newSheet.Cells.Item(31,3).Value = 43868.07152777778
$newSheet.Cells.Item(31,4).Value = 39
$newSheet.Cells.Item(31,5).Value = 0
$newSheet.Cells.Item(31,6).Value = 0
$newSheet.Cells.Item(32,2).Value = 'Singapore'
$newSheet.Cells.Item(32,3).Value = 43868.12708333333
$newSheet.Cells.Item(32,4).Value = 30
$newSheet.Cells.Item(32,5).Value = 0
$newSheet.Cells.Item(32,6).Value = 0
$newSheet.Cells.Item(33,1).Value = 'Hong Kong'
$newSheet.Cells.Item(33,2).Value = 'Hong Kong'
$newSheet.Cells.Item(33,3).Value = 43868.49513888889
$newSheet.Cells.Item(33,4).Value = 25
$newSheet.Cells.Item(33,5).Value = 1
$newSheet.Cells.Item(33,6).Value = 0
$newSheet.Cells.Item(34,2).Value = 'Japan'
$newSheet.Cells.Item(34,3).Value = 43868.17569444444
$newSheet.Cells.Item(34,4).Value = 25
$newSheet.Cells.Item(34,5).Value = 0
$newSheet.Cells.Item(34,6).Value = 1
$newSheet.Cells.Item(35,2).Value = 'Thailand'
$newSheet.Cells.Item(35,3).Value = 43865.64791666667
$newSheet.Cells.Item(35,4).Value = 25
$newSheet.Cells.Item(35,5).Value = 0
$newSheet.Cells.Item(35,6).Value = 5
$newSheet.Cells.Item(36,2).Value = 'South Korea'
$newSheet.Cells.Item(36,3).Value = 43868.12708333333
$newSheet.Cells.Item(36,4).Value = 24
$newSheet.Cells.Item(36,5).Value = 0
$newSheet.Cells.Item(36,6).Value = 1
$newSheet.Cells.Item(37,1).Value = 'Qinghai'
$newSheet.Cells.Item(37,2).Value = 'Mainland China'
$newSheet.Cells.Item(37,3).Value = 43867.092361111114
$newSheet.Cells.Item(37,4).Value = 18
$newSheet.Cells.Item(37,5).Value = 0
$newSheet.Cells.Item(37,6).Value = 3
$newSheet.Cells.Item(38,1).Value = 'Taiwan'
$newSheet.Cells.Item(38,2).Value = 'Taiwan'
$newSheet.Cells.Item(38,3).Value = 43867.62708333333
$newSheet.Cells.Item(38,4).Value = 16
$newSheet.Cells.Item(38,5).Value = 0
$newSheet.Cells.Item(38,6).Value = 1
$newSheet.Cells.Item(39,2).Value = 'Germany'
$newSheet.Cells.Item(39,3).Value = 43868.68958333333
$newSheet.Cells.Item(39,4).Value = 13
$newSheet.Cells.Item(39,5).Value = 0
$newSheet.Cells.Item(39,6).Value = 0
$newSheet.Cells.Item(40,2).Value = 'Malaysia'
$newSheet.Cells.Item(40,3).Value = 43868.38402777778
$newSheet.Cells.Item(40,4).Value = 12
$newSheet.Cells.Item(40,5).Value = 0
$newSheet.Cells.Item(40,6).Value = 1
$newSheet.Cells.Item(41,1).Value = 'Macau'
$newSheet.Cells.Item(41,2).Value = 'Macau'
$newSheet.Cells.Item(41,3).Value = 43867.59930555556
$newSheet.Cells.Item(41,4).Value = 10
$newSheet.Cells.Item(41,5).Value = 0
$newSheet.Cells.Item(41,6).Value = 1
$newSheet.Cells.Item(42,2).Value = 'Vietnam'
$newSheet.Cells.Item(42,3).Value = 43867.05069444444
$newSheet.Cells.Item(42,4).Value = 10
$newSheet.Cells.Item(42,5).Value = 0
$newSheet.Cells.Item(42,6).Value = 1
$newSheet.Cells.Item(43,2).Value = 'France'
$newSheet.Cells.Item(43,3).Value = 43862.07777777778
$newSheet.Cells.Item(43,4).Value = 6
$newSheet.Cells.Item(43,5).Value = 0
$newSheet.Cells.Item(43,6).Value = 0
$newSheet.Cells.Item(44,1).Value = 'Queensland'
$newSheet.Cells.Item(44,2).Value = 'Australia'
$newSheet.Cells.Item(44,3).Value = 43868.04375
$newSheet.Cells.Item(44,4).Value = 5
$newSheet.Cells.Item(44,5).Value = 0
$newSheet.Cells.Item(44,6).Value = 0
$newSheet.Cells.Item(45,2).Value = 'United Arab Emirates'
$newSheet.Cells.Item(45,3).Value = 43863.23819444444
$newSheet.Cells.Item(45,4).Value = 5
$newSheet.Cells.Item(45,5).Value = 0
$newSheet.Cells.Item(45,6).Value = 0
$newSheet.Cells.Item(46,1).Value = 'New South Wales'
$newSheet.Cells.Item(46,2).Value = 'Australia'
$newSheet.Cells.Item(46,3).Value = 43867.13402777778
$newSheet.Cells.Item(46,4).Value = 4
$newSheet.Cells.Item(46,5).Value = 0
$newSheet.Cells.Item(46,6).Value = 2
$newSheet.Cells.Item(47,1).Value = 'Victoria'
$newSheet.Cells.Item(47,2).Value = 'Australia'
$newSheet.Cells.Item(47,3).Value = 43862.75833333333
$newSheet.Cells.Item(47,4).Value = 4
$newSheet.Cells.Item(47,5).Value = 0
$newSheet.Cells.Item(47,6).Value = 0
$newSheet.Cells.Item(48,1).Value = 'British Columbia'
$newSheet.Cells.Item(48,2).Value = 'Canada'
$newSheet.Cells.Item(48,3).Value = 43868.23819444444
$newSheet.Cells.Item(48,4).Value = 4
$newSheet.Cells.Item(48,5).Value = 0
$newSheet.Cells.Item(48,6).Value = 0
$newSheet.Cells.Item(49,2).Value = 'India'
$newSheet.Cells.Item(49,3).Value = 43864.904861111114
$newSheet.Cells.Item(49,4).Value = 3
$newSheet.Cells.Item(49,5).Value = 0
$newSheet.Cells.Item(49,6).Value = 0
$newSheet.Cells.Item(50,2).Value = 'Italy'
$newSheet.Cells.Item(50,3).Value = 43868.74513888889
$newSheet.Cells.Item(50,4).Value = 3
$newSheet.Cells.Item(50,5).Value = 0
$newSheet.Cells.Item(50,6).Value = 0
$newSheet.Cells.Item(51,2).Value = 'Philippines'
$newSheet.Cells.Item(51,3).Value = 43868.69652777778
$newSheet.Cells.Item(51,4).Value = 3
$newSheet.Cells.Item(51,5).Value = 1
$newSheet.Cells.Item(51,6).Value = 0
$newSheet.Cells.Item(52,2).Value = 'UK'
$newSheet.Cells.Item(52,3).Value = 43868.77291666667
$newSheet.Cells.Item(52,4).Value = 3
$newSheet.Cells.Item(52,5).Value = 0
$newSheet.Cells.Item(52,6).Value = 0
$newSheet.Cells.Item(53,1).Value = 'South Australia'
$newSheet.Cells.Item(53,2).Value = 'Australia'
$newSheet.Cells.Item(53,3).Value = 43863.93958333333
$newSheet.Cells.Item(53,4).Value = 2
$newSheet.Cells.Item(53,5).Value = 0
$newSheet.Cells.Item(53,6).Value = 0
$newSheet.Cells.Item(54,1).Value = 'Toronto, ON'
$newSheet.Cells.Item(54,2).Value = 'Canada'
$newSheet.Cells.Item(54,3).Value = 43865.00902777778
$newSheet.Cells.Item(54,4).Value = 2
$newSheet.Cells.Item(54,5).Value = 0
$newSheet.Cells.Item(54,6).Value = 0
$newSheet.Cells.Item(55,2).Value = 'Russia'
$newSheet.Cells.Item(55,3).Value = 43861.67569444444
$newSheet.Cells.Item(55,4).Value = 2
$newSheet.Cells.Item(55,5).Value = 0
$newSheet.Cells.Item(55,6).Value = 0
$newSheet.Cells.Item(56,1).Value = 'Chicago, IL'
$newSheet.Cells.Item(56,2).Value = 'US'
$newSheet.Cells.Item(56,3).Value = 43862.82152777778
$newSheet.Cells.Item(56,4).Value = 2
$newSheet.Cells.Item(56,5).Value = 0
$newSheet.Cells.Item(56,6).Value = 0
$newSheet.Cells.Item(57,1).Value = 'San Benito, CA'
$newSheet.Cells.Item(57,2).Value = 'US'
$newSheet.Cells.Item(57,3).Value = 43864.16180555556
$newSheet.Cells.Item(57,4).Value = 2
$newSheet.Cells.Item(57,5).Value = 0
$newSheet.Cells.Item(57,6).Value = 0
$newSheet.Cells.Item(58,1).Value = 'Santa Clara, CA'
$newSheet.Cells.Item(58,2).Value = 'US'
$newSheet.Cells.Item(58,3).Value = 43864.029861111114
$newSheet.Cells.Item(58,4).Value = 2
$newSheet.Cells.Item(58,5).Value = 0
$newSheet.Cells.Item(58,6).Value = 0
$newSheet.Cells.Item(59,2).Value = 'Belgium'
$newSheet.Cells.Item(59,3).Value = 43865.654861111114
$newSheet.Cells.Item(59,4).Value = 1
$newSheet.Cells.Item(59,5).Value = 0
$newSheet.Cells.Item(59,6).Value = 0
$newSheet.Cells.Item(60,2).Value = 'Cambodia'
$newSheet.Cells.Item(60,3).Value = 43861.34375
$newSheet.Cells.Item(60,4).Value = 1
$newSheet.Cells.Item(60,5).Value = 0
$newSheet.Cells.Item(60,6).Value = 0
$newSheet.Cells.Item(61,1).Value = 'London, ON'
$newSheet.Cells.Item(61,2).Value = 'Canada'
$newSheet.Cells.Item(61,3).Value = 43865.00208333333
$newSheet.Cells.Item(61,4).Value = 1
$newSheet.Cells.Item(61,5).Value = 0
$newSheet.Cells.Item(61,6).Value = 0
$newSheet.Cells.Item(62,2).Value = 'Finland'
$newSheet.Cells.Item(62,3).Value = 43861.34375
$newSheet.Cells.Item(62,4).Value = 1
$newSheet.Cells.Item(62,5).Value = 0
$newSheet.Cells.Item(62,6).Value = 0
$newSheet.Cells.Item(63,1).Value = 'Tibet'
$newSheet.Cells.Item(63,2).Value = 'Mainland China'
$newSheet.Cells.Item(63,3).Value = 43862.07777777778
$newSheet.Cells.Item(63,4).Value = 1
$newSheet.Cells.Item(63,5).Value = 0
$newSheet.Cells.Item(63,6).Value = 0
$newSheet.Cells.Item(64,2).Value = 'Nepal'
$newSheet.Cells.Item(64,3).Value = 43861.34375
$newSheet.Cells.Item(64,4).Value = 1
$newSheet.Cells.Item(64,5).Value = 0
$newSheet.Cells.Item(64,6).Value = 0
$newSheet.Cells.Item(65,2).Value = 'Spain'
$newSheet.Cells.Item(65,3).Value = 43862.98819444444
$newSheet.Cells.Item(65,4).Value = 1
$newSheet.Cells.Item(65,5).Value = 0
$newSheet.Cells.Item(65,6).Value = 0
$newSheet.Cells.Item(66,2).Value = 'Sri Lanka'
$newSheet.Cells.Item(66,3).Value = 43861.34375
$newSheet.Cells.Item(66,4).Value = 1
$newSheet.Cells.Item(66,5).Value = 0
$newSheet.Cells.Item(66,6).Value = 0
$newSheet.Cells.Item(67,2).Value = 'Sweden'
$newSheet.Cells.Item(67,3).Value = 43862.092361111114
$newSheet.Cells.Item(67,4).Value = 1
$newSheet.Cells.Item(67,5).Value = 0
$newSheet.Cells.Item(67,6).Value = 0
$newSheet.Cells.Item(68,1).Value = 'Boston, MA'
$newSheet.Cells.Item(68,2).Value = 'US'
$newSheet.Cells.Item(68,3).Value = 43862.82152777778
$newSheet.Cells.Item(68,4).Value = 1
$newSheet.Cells.Item(68,5).Value = 0
$newSheet.Cells.Item(68,6).Value = 0
$newSheet.Cells.Item(69,1).Value = 'Los Angeles, CA'
$newSheet.Cells.Item(69,2).Value = 'US'
$newSheet.Cells.Item(69,3).Value = 43862.82847222222
$newSheet.Cells.Item(69,4).Value = 1
$newSheet.Cells.Item(69,5).Value = 0
$newSheet.Cells.Item(69,6).Value = 0
$newSheet.Cells.Item(70,1).Value = 'Madison, WI'
$newSheet.Cells.Item(70,2).Value = 'US'
$newSheet.Cells.Item(70,3).Value = 43866.91180555556
$newSheet.Cells.Item(70,4).Value = 1
$newSheet.Cells.Item(70,5).Value = 0
$newSheet.Cells.Item(70,6).Value = 0
$newSheet.Cells.Item(71,1).Value = 'Orange, CA'
$newSheet.Cells.Item(71,2).Value = 'US'
$newSheet.Cells.Item(71,3).Value = 43862.82847222222
$newSheet.Cells.Item(71,4).Value = 1
$newSheet.Cells.Item(71,5).Value = 0
$newSheet.Cells.Item(71,6).Value = 0
$newSheet.Cells.Item(72,1).Value = 'Seattle, WA'
$newSheet.Cells.Item(72,2).Value = 'US'
$newSheet.Cells.Item(72,3).Value = 43862.82152777778
$newSheet.Cells.Item(72,4).Value = 1
$newSheet.Cells.Item(72,5).Value = 0
$newSheet.Cells.Item(72,6).Value = 0
$newSheet.Cells.Item(73,1).Value = 'Tempe, AZ'
$newSheet.Cells.Item(73,2).Value = 'US'
$newSheet.Cells.Item(73,3).Value = 43862.82152777778
$newSheet.Cells.Item(73,4).Value = 1
$newSheet.Cells.Item(73,5).Value = 0
$newSheet.Cells.Item(73,6).Value = 0

# Row 73 is brand new (sheet grew from 72 to 73 data rows) - give its date cell the same
# date/time number format used by the rest of column C
$newSheet.Cells.Item(73,3).NumberFormat = "m/d/yy h:mm"

# 3) Match the widened columns captured for the new sheet
$newSheet.Columns.Item(1).ColumnWidth = 22.276041666666668
$newSheet.Columns.Item(2).ColumnWidth = 25.608072916666668
$newSheet.Columns.Item(3).ColumnWidth = 17.385416666666668
$newSheet.Columns.Item(4).ColumnWidth = 15.721354166666666
$newSheet.Columns.Item(5).ColumnWidth = 13.053385416666666
$newSheet.Columns.Item(6).ColumnWidth = 12.385416666666666

# 4) Put the view/selection on the new (now active) sheet where it was left
$newSheet.Activate()
$newSheet.Cells.Item(87,13).Select()

# 5) The previous "latest" sheet is no longer the active tab - park its selection
#    on the header row, matching the saved view
$srcSheet.Activate()
$srcSheet.Range("A1:F1").Select()

# 6) Leave the new sheet active/selected as the final state
$newSheet.Activate()
